$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8289269051321928
$ws.Range("B4").Value = 0.8419243986254296
$ws.Range("B5").Value = 0.3185185185185185
$ws.Range("B6").Value = 0.9645669291338582
$ws.Range("B7").Value = 0.872852233676976
